$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.136.98'
$ws.Range("E2").Value = '  +1.53%  '
$ws.Range("D3").Value = '3.771.22'
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '622.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.75%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '164.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '3.767.83'
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E10").Value = '  +1.24%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.450'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '0.0000247'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("E14").Value = '  -0.64%  '
$ws.Range("D15").Value = '4.407.59'
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("D16").Value = '3.771.15'
$ws.Range("E16").Value = '  -1.34%  '
$ws.Range("D17").Value = '69.170.49'
$ws.Range("E18").Value = '  -3.34%  '
$ws.Range("E19").Value = '  +0.43%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '467.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.52%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '9.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '0.699'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  +2.37%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '83.06'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '12.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("E27").Value = '  +2.40%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '10.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("D30").Value = '3.919.73'
$ws.Range("E30").Value = '  -0.49%  '
$ws.Range("E31").Value = '  +0.44%  '
$ws.Range("E32").Value = '  +0.75%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '7.27'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '28.83'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = '3.724.61'
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '8.95'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.39%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.156'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.60%  '
$ws.Range("E39").Value = '  +2.22%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '3.34'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.57%  '
$ws.Range("E41").Value = '  -0.78%  '
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.964'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.20%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '0.299'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '153.78'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '42.83'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.38%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '46.63'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.65%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '1.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.09%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '8.38'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.60%  '
$ws.Range("E51").Value = '  +0.68%  '
